$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlLeft   = -4131
$xlCenter = -4108

# --- Insert a new row above row 1, and a new column before column D ---
# Before: row1 = A1(會計日期) C1(blank,style "center/center") ; row2 = A2:G2 header labels
# After insert: row1 is blank/new, old row1 -> row2, old row2 -> row3
#               old column D (index4) -> E, etc. A new blank column D is created
$ws.Rows.Item(1).Insert()
$ws.Columns.Item(4).Insert()

# --- Row 1 ---
$ws.Range("A1").Value = "會計日期"
$ws.Range("A1").HorizontalAlignment = $xlLeft
$ws.Range("A1").VerticalAlignment = $xlCenter

$ws.Range("B1").ClearContents()

$ws.Range("C1").Value = ""
$ws.Range("C1").HorizontalAlignment = $xlCenter
$ws.Range("C1").VerticalAlignment = $xlCenter

$ws.Range("D1").Value = ""
$ws.Range("D1").HorizontalAlignment = $xlCenter
$ws.Range("D1").VerticalAlignment = $xlCenter

# --- Row 2 (was old row1, now shifted down & right) ---
$ws.Range("A2").Clear()

$ws.Range("C2").Value = ""
$ws.Range("C2").HorizontalAlignment = $xlCenter
$ws.Range("C2").VerticalAlignment = $xlCenter

$ws.Range("D2").Value = ""
$ws.Range("D2").HorizontalAlignment = $xlCenter
$ws.Range("D2").VerticalAlignment = $xlCenter

$ws.Range("E2").Value = "今日"
$ws.Range("F2").Value = "今日"
$ws.Range("G2").Value = "昨日"
$ws.Range("H2").Value = "昨日"
$ws.Range("E2:H2").HorizontalAlignment = $xlLeft
$ws.Range("E2:H2").VerticalAlignment = $xlCenter

# --- Row 3 (was old row2, now shifted right from column D onward) ---
$ws.Range("A3").Value = "會計備份日"
$ws.Range("B3").Value = "戶號"
$ws.Range("C3").Value = "戶名"
$ws.Range("D3").Value = "額度"
$ws.Range("E3").Value = "暫收款支票"
$ws.Range("F3").Value = "暫收款非支票"
$ws.Range("G3").Value = "暫收款支票"
$ws.Range("H3").Value = "暫收款非支票"
$ws.Range("A3:H3").HorizontalAlignment = $xlLeft
$ws.Range("A3:H3").VerticalAlignment = $xlCenter

# --- Column widths (new column D matches its neighbour C: default width 9) ---
$ws.Columns.Item(4).ColumnWidth = 8.2875

# --- Selection, matching the state recorded after the edit ---
$ws.Range("H4").Select()
